# Append the latest daily rows (45971-45978) to each of the 4 sheets and
# fill in the previously-placeholder value for row 105 (date 45968), which
# had been recorded as 0 pending data.

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Data per sheet: index 1..4 matching Worksheets.Item(1..4)
# Each entry: row105 fix-up value, then new rows 106-111 as (date, value)
$sheetsData = @(
    @{
        Row105 = 638
        NewRows = @(
            @(45971, 521),
            @(45972, 521),
            @(45973, 521),
            @(45974, 518),
            @(45975, 516),
            @(45978, 513)
        )
    },
    @{
        Row105 = 7260
        NewRows = @(
            @(45971, 7339),
            @(45972, 7121),
            @(45973, 7175),
            @(45974, 7214),
            @(45975, 7134),
            @(45978, 7081)
        )
    },
    @{
        Row105 = 2781
        NewRows = @(
            @(45971, 2917),
            @(45972, 2926),
            @(45973, 2922),
            @(45974, 2975),
            @(45975, 2930),
            @(45978, 2867)
        )
    },
    @{
        Row105 = 1282
        NewRows = @(
            @(45971, 1238),
            @(45972, 1223),
            @(45973, 1218),
            @(45974, 1243),
            @(45975, 1260),
            @(45978, 1333)
        )
    }
)

for ($i = 0; $i -lt 4; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $sheetInfo = $sheetsData[$i]

    # Fix up row 105, column B (was a 0 placeholder, now the real value)
    $ws.Cells.Item(105, 2).Value = $sheetInfo.Row105

    # Append rows 106-111
    $r = 106
    foreach ($pair in $sheetInfo.NewRows) {
        $ws.Cells.Item($r, 1).Value = $pair[0]
        $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
        $ws.Cells.Item($r, 2).Value = $pair[1]
        $r = $r + 1
    }
}
